# Update workbook "Översikt FILIPSTAD" worksheet with the latest scrape data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") for existing data rows 2-455 moves from 45190 to 45192.
$ws.Range("C2:C455").Value = 45192

# 2) Row 455 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(455).RowHeight = 15

# 3) A new data row (456) is appended for case "A 44792-2023".
$ws.Cells.Item(456, 1).Value = "A 44792-2023"

$ws.Cells.Item(456, 2).Value = 45190
$ws.Cells.Item(456, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(456, 3).Value = 45192
$ws.Cells.Item(456, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(456, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item(456, 5).Value = "FILIPSTAD"

$ws.Cells.Item(456, 7).Value = 18.4

$ws.Cells.Item(456, 8).Value = 0
$ws.Cells.Item(456, 9).Value = 0
$ws.Cells.Item(456, 10).Value = 0
$ws.Cells.Item(456, 11).Value = 0
$ws.Cells.Item(456, 12).Value = 0
$ws.Cells.Item(456, 13).Value = 0
$ws.Cells.Item(456, 14).Value = 0
$ws.Cells.Item(456, 15).Value = 0
$ws.Cells.Item(456, 16).Value = 0
$ws.Cells.Item(456, 17).Value = 0

$ws.Cells.Item(456, 18).WrapText = $true
$ws.Cells.Item(456, 18).Value = ""
